$d = $word.ActiveDocument

# The paragraph currently reads "PRUEBA 1.1" where the last run holds ".1".
# We need it to read "PRUEBA 1.2", with the final "1" becoming "2" and that
# "2" living in its own run (splitting the original ".1" run into "." and "2").
$text = $d.Content.Text
$idx = $text.LastIndexOf("1")
$r = $d.Range($idx, $idx + 1)

# Temporarily toggle a character attribute so the replaced character is
# forced into its own run instead of being re-merged with its neighbor.
$r.Bold = 1
$r.Text = "2"

$r2 = $d.Range($idx, $idx + 1)
$r2.Bold = 0
